# feat: add 2022-Q3 data
#
# The workbook originally has two sheets: "总计" (summary) and "2022-Q1"
# (the Q1 fund-holding detail sheet). This edit inserts a brand-new
# "2022-Q3" detail sheet between them (so the tab order becomes
# 总计 / 2022-Q3 / 2022-Q1) and records the Q3 summary figures on the
# "总计" sheet as the new row 2, pushing the old Q1 summary row down to
# row 3.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # 总计
$wsQ1 = $wb.Worksheets.Item(2)     # 2022-Q1 (untouched detail sheet)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" detail sheet right after "总计" so the
#    final tab order is 总计, 2022-Q3, 2022-Q1.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $ws1)
$q3.Name = "2022-Q3"

# Match the page margins used across the rest of the workbook (0.75in /
# 1in / 0.5in) instead of the engine's blank-sheet defaults.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Carry over the bold/centered/bordered header style and the "index"
# column style from the "总计" sheet so the new sheet matches the
# workbook's existing look instead of picking up a brand-new style.
$ws1.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$q3.Range("A2:A3").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Fund figures are stored as text (as the rest of the workbook does),
# so force a text number format before writing the numeric-looking
# values.
$q3.Range("B2:G3").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "159869"
$q3.Range("C2").Value = "华夏中证动漫游戏ETF"
$q3.Range("D2").Value = "6.35"
$q3.Range("E2").Value = "99.31"
$q3.Range("F2").Value = "3.11"
$q3.Range("G2").Value = "0.1975"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "516010"
$q3.Range("C3").Value = "国泰中证动漫游戏ETF"
$q3.Range("D3").Value = "3.78"
$q3.Range("E3").Value = "97.86"
$q3.Range("F3").Value = "3.05"
$q3.Range("G3").Value = "0.1153"
$q3.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push the existing 2022-Q1 row
#    down to row 3, and write the new 2022-Q3 totals into row 2.
# ---------------------------------------------------------------------
$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q1"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.09

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.31

# ---------------------------------------------------------------------
# 3. Restore "总计" as the active tab (it was unaffected by the edit).
# ---------------------------------------------------------------------
$ws1.Activate()
